# Auto-generated: refresh market-data-driven columns (H-N) across all profession sheets
# Source: scheduled runner that pulls updated Universalis pricing data
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2717.1667
$ws.Range("J32").Value = 2771.8
$ws.Range("L32").Value = 2771.8
$ws.Range("N32").Value = -3423.8
$ws.Range("H38").Value = 25291.143
$ws.Range("I38").Value = 33389.9
$ws.Range("K38").Value = 100169.7
$ws.Range("M38").Value = -99797.70000000001
$ws.Range("H53").Value = 894
$ws.Range("J53").Value = 922.6667
$ws.Range("L53").Value = 922.6667
$ws.Range("N53").Value = -2196.6667
$ws.Range("H88").Value = 4179.5
$ws.Range("J88").Value = 4598.6
$ws.Range("L88").Value = 4598.6
$ws.Range("N88").Value = -5410.6
$ws.Range("H91").Value = 4179.5
$ws.Range("J91").Value = 4598.6
$ws.Range("L91").Value = 4598.6
$ws.Range("N91").Value = -7406.6
$ws.Range("H98").Value = 1361.0286
$ws.Range("I98").Value = 1201.8572
$ws.Range("K98").Value = 1201.8572
$ws.Range("M98").Value = 296.1428000000001
$ws.Range("H113").Value = 8251.25
$ws.Range("I113").Value = 7499.5
$ws.Range("J113").Value = 9003
$ws.Range("K113").Value = 7499.5
$ws.Range("L113").Value = 9003
$ws.Range("M113").Value = -4245.5
$ws.Range("N113").Value = -15511
$ws.Range("H116").Value = 13441.823
$ws.Range("I116").Value = 7533
$ws.Range("K116").Value = 7533
$ws.Range("M116").Value = -4091
$ws.Range("H122").Value = 1361.0286
$ws.Range("I122").Value = 1201.8572
$ws.Range("K122").Value = 3605.5716
$ws.Range("M122").Value = -1155.5716
$ws.Range("H132").Value = 4792.7188
$ws.Range("I132").Value = 2339.682
$ws.Range("J132").Value = 10189.4
$ws.Range("K132").Value = 7019.045999999999
$ws.Range("L132").Value = 30568.2
$ws.Range("M132").Value = -4489.045999999999
$ws.Range("N132").Value = -35628.2
$ws.Range("H138").Value = 11065.363
$ws.Range("I138").Value = 3059.4
$ws.Range("J138").Value = 17737
$ws.Range("K138").Value = 9178.200000000001
$ws.Range("L138").Value = 53211
$ws.Range("M138").Value = -4038.200000000001
$ws.Range("N138").Value = -63491

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2395
$ws.Range("I45").Value = 1243.75
$ws.Range("J45").Value = 7000
$ws.Range("K45").Value = 1243.75
$ws.Range("L45").Value = 7000
$ws.Range("M45").Value = -866.75
$ws.Range("N45").Value = -7754
$ws.Range("H61").Value = 50716710
$ws.Range("I61").Value = 77780720
$ws.Range("J61").Value = 2001494.2
$ws.Range("K61").Value = 77780720
$ws.Range("L61").Value = 2001494.2
$ws.Range("M61").Value = -77780508
$ws.Range("N61").Value = -2001918.2
$ws.Range("H102").Value = 2268.6667
$ws.Range("J102").Value = 3193.8572
$ws.Range("L102").Value = 3193.8572
$ws.Range("N102").Value = -6437.8572
$ws.Range("H122").Value = 2310
$ws.Range("I122").Value = 2255.0625
$ws.Range("K122").Value = 6765.1875
$ws.Range("M122").Value = -4315.1875
$ws.Range("H132").Value = 2086860.8
$ws.Range("I132").Value = 3834.2927
$ws.Range("J132").Value = 14287444
$ws.Range("K132").Value = 11502.8781
$ws.Range("L132").Value = 42862332
$ws.Range("M132").Value = -8972.8781
$ws.Range("N132").Value = -42867392
$ws.Range("H136").Value = 50716710
$ws.Range("I136").Value = 77780720
$ws.Range("J136").Value = 2001494.2
$ws.Range("K136").Value = 233342160
$ws.Range("L136").Value = 6004482.6
$ws.Range("M136").Value = -233339610
$ws.Range("N136").Value = -6009582.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1433267.9
$ws.Range("I105").Value = 2079471.5
$ws.Range("J105").Value = 11619.8
$ws.Range("K105").Value = 2079471.5
$ws.Range("L105").Value = 11619.8
$ws.Range("M105").Value = -2077724.5
$ws.Range("N105").Value = -15113.8
$ws.Range("H134").Value = 3848363.8
$ws.Range("I134").Value = 1691.375
$ws.Range("K134").Value = 5074.125
$ws.Range("M134").Value = -2539.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 14289851
$ws.Range("I16").Value = 14289851
$ws.Range("K16").Value = 14289851
$ws.Range("M16").Value = -14289564
$ws.Range("H31").Value = 33673364
$ws.Range("I31").Value = 40003316
$ws.Range("J31").Value = 2023603.8
$ws.Range("K31").Value = 40003316
$ws.Range("L31").Value = 2023603.8
$ws.Range("M31").Value = -40003021
$ws.Range("N31").Value = -2024193.8
$ws.Range("H34").Value = 33673364
$ws.Range("I34").Value = 40003316
$ws.Range("J34").Value = 2023603.8
$ws.Range("K34").Value = 40003316
$ws.Range("L34").Value = 2023603.8
$ws.Range("M34").Value = -40003114
$ws.Range("N34").Value = -2024007.8
$ws.Range("H51").Value = 27989.5
$ws.Range("I51").Value = 27989.5
$ws.Range("K51").Value = 27989.5
$ws.Range("M51").Value = -27253.5
$ws.Range("H61").Value = 27989.5
$ws.Range("I61").Value = 27989.5
$ws.Range("K61").Value = 27989.5
$ws.Range("M61").Value = -27641.5
$ws.Range("H102").Value = 83490.5
$ws.Range("J102").Value = 83490.5
$ws.Range("L102").Value = 83490.5
$ws.Range("N102").Value = -88358.5
$ws.Range("H107").Value = 1686.8125
$ws.Range("I107").Value = 1493.0541
$ws.Range("J107").Value = 2338.5454
$ws.Range("K107").Value = 1493.0541
$ws.Range("L107").Value = 2338.5454
$ws.Range("M107").Value = 426.9458999999999
$ws.Range("N107").Value = -6178.5454
$ws.Range("H113").Value = 14289851
$ws.Range("I113").Value = 14289851
$ws.Range("K113").Value = 14289851
$ws.Range("M113").Value = -14287681
$ws.Range("H122").Value = 2569.8
$ws.Range("I122").Value = 2378.682
$ws.Range("K122").Value = 7136.045999999999
$ws.Range("M122").Value = -4686.045999999999
$ws.Range("H132").Value = 2157.8
$ws.Range("I132").Value = 1772.25
$ws.Range("K132").Value = 5316.75
$ws.Range("M132").Value = -2786.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 4400.4
$ws.Range("J68").Value = 4750
$ws.Range("L68").Value = 14250
$ws.Range("N68").Value = -15872
$ws.Range("H71").Value = 4400.4
$ws.Range("J71").Value = 4750
$ws.Range("L71").Value = 42750
$ws.Range("N71").Value = -50862
$ws.Range("H121").Value = 4522.5557
$ws.Range("J121").Value = 5294.7334
$ws.Range("L121").Value = 15884.2002
$ws.Range("N121").Value = -18504.2002
$ws.Range("H141").Value = 6365.4546
$ws.Range("I141").Value = 3668.7
$ws.Range("K141").Value = 11006.1
$ws.Range("M141").Value = -5826.099999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 26089238
$ws.Range("I80").Value = 2033.5834
$ws.Range("J80").Value = 54548004
$ws.Range("K80").Value = 2033.5834
$ws.Range("L80").Value = 54548004
$ws.Range("M80").Value = -1035.5834
$ws.Range("N80").Value = -54550000
$ws.Range("H83").Value = 26089238
$ws.Range("I83").Value = 2033.5834
$ws.Range("J83").Value = 54548004
$ws.Range("K83").Value = 10167.917
$ws.Range("L83").Value = 272740020
$ws.Range("M83").Value = -5175.916999999999
$ws.Range("N83").Value = -272750004
$ws.Range("H107").Value = 1824.125
$ws.Range("I107").Value = 1799
$ws.Range("K107").Value = 1799
$ws.Range("M107").Value = 121
$ws.Range("H122").Value = 6662.6665
$ws.Range("I122").Value = 7494
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 22482
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -20032
$ws.Range("N122").Value = -19900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3245.375
$ws.Range("J46").Value = 3498.5
$ws.Range("L46").Value = 3498.5
$ws.Range("N46").Value = -3874.5
$ws.Range("H68").Value = 3789802
$ws.Range("I68").Value = 4168381.8
$ws.Range("J68").Value = 4003
$ws.Range("K68").Value = 4168381.8
$ws.Range("L68").Value = 4003
$ws.Range("M68").Value = -4167632.8
$ws.Range("N68").Value = -5501
$ws.Range("H71").Value = 3789802
$ws.Range("I71").Value = 4168381.8
$ws.Range("J71").Value = 4003
$ws.Range("K71").Value = 20841909
$ws.Range("L71").Value = 20015
$ws.Range("M71").Value = -20838165
$ws.Range("N71").Value = -27503
$ws.Range("H93").Value = 3273294.5
$ws.Range("J93").Value = 6180558
$ws.Range("L93").Value = 6180558
$ws.Range("N93").Value = -6183054
$ws.Range("H132").Value = 3800.2222
$ws.Range("I132").Value = 2279.4167
$ws.Range("K132").Value = 6838.250100000001
$ws.Range("M132").Value = -4308.250100000001
$ws.Range("H136").Value = 5692.25
$ws.Range("I136").Value = 5400.7
$ws.Range("J136").Value = 7150
$ws.Range("K136").Value = 16202.1
$ws.Range("L136").Value = 21450
$ws.Range("M136").Value = -13652.1
$ws.Range("N136").Value = -26550

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 28796.8
$ws.Range("J41").Value = 28497.75
$ws.Range("L41").Value = 28497.75
$ws.Range("N41").Value = -29277.75
$ws.Range("H51").Value = 61475
$ws.Range("J51").Value = 61475
$ws.Range("L51").Value = 61475
$ws.Range("N51").Value = -62495
$ws.Range("H81").Value = 1469
$ws.Range("I81").Value = 1469
$ws.Range("K81").Value = 2938
$ws.Range("M81").Value = -1877
$ws.Range("H84").Value = 1469
$ws.Range("I84").Value = 1469
$ws.Range("K84").Value = 14690
$ws.Range("M84").Value = -9386
$ws.Range("H100").Value = 794458.8
$ws.Range("I100").Value = 1672.4166
$ws.Range("J100").Value = 5551177
$ws.Range("K100").Value = 3344.8332
$ws.Range("L100").Value = 11102354
$ws.Range("M100").Value = -2803.8332
$ws.Range("N100").Value = -11103436
$ws.Range("H107").Value = 2839.6667
$ws.Range("I107").Value = 1871
$ws.Range("J107").Value = 3254.8096
$ws.Range("K107").Value = 5613
$ws.Range("L107").Value = 9764.4288
$ws.Range("M107").Value = -3693
$ws.Range("N107").Value = -13604.4288

